$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: D1 label changed from "Cross Validation Mean Accuracy" to "Mean Cross Validation Accuracy"
$ws.Range("D1").Value = "Mean Cross Validation Accuracy"

# Row 8 used to be a lone "KNeighbours Classifier (Count Vectorizer + TfidfTransformer)" row with
# no metrics. It becomes the "KNeighbours Classifier (Count Vectorizer)" row with full metrics
# (using the same number formatting as the rest of the data rows), and the model names in the
# rows below shuffle around as results for the remaining KNeighbours / Gradient Boosting /
# Xgboost model variants are filled in.
$ws.Range("A8").Value = "KNeighbours Classifier (Count Vectorizer)"
$ws.Range("B2").Copy()
$ws.Range("B8:D8").PasteSpecial(-4122)
$ws.Range("B8").Value = 0.91916624300965899
$ws.Range("C8").Value = 0.94457157386219104
$ws.Range("D8").Value = 0.83082393952284905

$ws.Range("A9").Value = "KNeighbours Classifier (Count Vectorizer + TfidfTransformer)"

$ws.Range("A10").Value = "Gradient Boosting Classifier (Count Vectorizer)"
$ws.Range("C2").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = 0.90318344471214596

$ws.Range("A11").Value = "Gradient Boosting Classifier (Count Vectorizer + TfidfTransformer)"

$ws.Range("A12").Value = "Xgboost Classifier (Count Vectorizer)"
$ws.Range("C2").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = 0.93089596444564304

$ws.Range("A13").Value = "Xgboost Classifier (Count Vectorizer + TfidfTransformer)"
$ws.Range("C2").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = 0.92825155622047995

# Selection moved to D11 (last cell touched interactively before save)
$ws.Range("D11").Select()
